$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Hunk 1: "Corso di " + "Fondamenti di intelligenza artificiale" -> one run
# ---------------------------------------------------------------------------
$t1 = "Corso di Fondamenti di intelligenza artificiale"
$null = $d.Content.Find.Execute($t1, $true, $false, $false, $false, $false, $true, 1, $false, $t1, 2)

# ---------------------------------------------------------------------------
# Hunk 3: "Etichettatore della bottiglia" + " nell'immagine" -> one run
# ---------------------------------------------------------------------------
$t3 = "Etichettatore della bottiglia nell" + [char]0x2019 + "immagine"
$null = $d.Content.Find.Execute($t3, $true, $false, $false, $false, $false, $true, 1, $false, $t3, 2)

# ---------------------------------------------------------------------------
# Hunk 4: the five "sottocartelle" runs -> one run
# ---------------------------------------------------------------------------
$t4 = [char]0x201D + " ovvero, contiene due sottocartelle: " + [char]0x201C + "bianco" + [char]0x201D + " e " + [char]0x201C + "rosso" + [char]0x201D + " contenenti rispettivamente 20 immagini di bottiglie di vini appartenenti alle categorie descritte dal nome della cartella."
$null = $d.Content.Find.Execute($t4, $true, $false, $false, $false, $false, $true, 1, $false, $t4, 2)

# ---------------------------------------------------------------------------
# Hunk 6: " che potrebbe rivelarsi utile anche " + "ad altri utenti." -> one
# run, and drop the _GoBack bookmark that used to sit between them (it will
# be re-created at its new home in hunk 2 below).
# ---------------------------------------------------------------------------
$oldGoBack = $d.Bookmarks.Item("_GoBack")
$oldGoBack.Delete()

$t6 = " che potrebbe rivelarsi utile anche ad altri utenti."
$null = $d.Content.Find.Execute($t6, $true, $false, $false, $false, $false, $true, 1, $false, $t6, 2)

# ---------------------------------------------------------------------------
# Hunk 2: fix the "`e" -> "è" typo and re-split the run into three pieces
# with the _GoBack bookmark sitting between piece 2 and piece 3.
# ---------------------------------------------------------------------------
$searchOld = "Nel caso di riconoscimento specifico, la scopo " + [char]0x60 + "e qu"
$r = $d.Content.Duplicate
$null = $r.Find.Execute($searchOld, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$base = $r.Start
$r.Text = "Nel caso di riconoscimento specifico, la scopo " + [char]0x00E8 + " qu"

# Character offsets (relative to $base) of the run boundaries we need to
# recreate: "Nel caso di ricon" | "oscimento specifico, la scopo è" |
# (bookmark) | " qu" | "ello di identificare un istanza " |
# "di un particolare oggetto, persona, luogo" | "."
$p1 = $base + 17
$p2 = $base + 48
$p3 = $base + 51
$p4 = $base + 83
$p5 = $base + 124

$rr1 = $d.Range($p1, $p1)
$d.Bookmarks.Add("TempSplit1", $rr1)
$d.Bookmarks.Item("TempSplit1").Delete()

$rr2 = $d.Range($p2, $p2)
$d.Bookmarks.Add("_GoBack", $rr2)

$rr3 = $d.Range($p3, $p3)
$d.Bookmarks.Add("TempSplit3", $rr3)
$d.Bookmarks.Item("TempSplit3").Delete()

$rr4 = $d.Range($p4, $p4)
$d.Bookmarks.Add("TempSplit4", $rr4)
$d.Bookmarks.Item("TempSplit4").Delete()

$rr5 = $d.Range($p5, $p5)
$d.Bookmarks.Add("TempSplit5", $rr5)
$d.Bookmarks.Item("TempSplit5").Delete()

# ---------------------------------------------------------------------------
# Hunk 5a: first drawing (rId7 / "Immagine 3") gets <w:noProof/>
# ---------------------------------------------------------------------------
$shp1 = $d.InlineShapes.Item(2)
$shp1.Range.NoProofing = 1

# ---------------------------------------------------------------------------
# Hunk 5b: second drawing (rId8 / "Immagine 4") gets <w:noProof/> and
# <w:lang w:eastAsia="it-IT"/>
# ---------------------------------------------------------------------------
$shp2 = $d.InlineShapes.Item(3)
$shp2.Range.NoProofing = 1
$shp2.Range.LanguageIDFarEast = "it-IT"

# ---------------------------------------------------------------------------
# styles.xml: Carpredefinitoparagrafo (Default Paragraph Font) gains
# <w:semiHidden/>
# ---------------------------------------------------------------------------
$s = $d.Styles.Item("Carpredefinitoparagrafo")
try {
  $s.Hidden = $true
} catch {
}

Write-Host "done"
